$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D, E, F) before the existing "Terms Typically Offered"
# column, shifting it from D to G.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = 'Corequisites'
$ws.Range("E1").Value = 'Concurrent'
$ws.Range("F1").Value = 'Recommended'
$ws.Range("G1").Value = 'Terms Typically Offered'

# Data rows: split the Corequisite/Concurrent clauses that had been embedded in the
# Prerequisites (C) / Terms Typically Offered (old D, now G) text into their own columns.
# Row 2
$ws.Range("D2").Value = 'NA'
$ws.Range("E2").Value = 'NA'
$ws.Range("F2").Value = 'NA'
$ws.Range("G2").Value = 'F'

# Row 3
$ws.Range("D3").Value = 'NA'
$ws.Range("E3").Value = 'NA'
$ws.Range("F3").Value = 'NA'
$ws.Range("G3").Value = 'F, SP'

# Row 4
$ws.Range("D4").Value = 'NA'
$ws.Range("E4").Value = 'NA'
$ws.Range("F4").Value = 'NA'
$ws.Range("G4").Value = 'W'

# Row 5
$ws.Range("D5").Value = 'LA 170.'
$ws.Range("E5").Value = 'NA'
$ws.Range("F5").Value = 'NA'
$ws.Range("G5").Value = 'F'

# Row 6
$ws.Range("C6").Value = 'LA 202, LA 243.'
$ws.Range("D6").Value = 'NA'
$ws.Range("E6").Value = 'LA 241.'
$ws.Range("F6").Value = 'NA'
$ws.Range("G6").Value = 'W '

# Row 7
$ws.Range("C7").Value = 'LA 203, LA 241.'
$ws.Range("D7").Value = 'NA'
$ws.Range("E7").Value = 'LA 242.'
$ws.Range("F7").Value = 'NA'
$ws.Range("G7").Value = 'SP '

# Row 8
$ws.Range("D8").Value = 'NA'
$ws.Range("E8").Value = 'NA'
$ws.Range("F8").Value = 'NA'
$ws.Range("G8").Value = 'W, SU'

# Row 9
$ws.Range("D9").Value = 'NA'
$ws.Range("E9").Value = 'NA'
$ws.Range("F9").Value = 'NA'
$ws.Range("G9").Value = 'F, SP, SU'

# Row 10
$ws.Range("D10").Value = 'NA'
$ws.Range("E10").Value = 'NA'
$ws.Range("F10").Value = 'NA'
$ws.Range("G10").Value = 'TBD'

# Row 11
$ws.Range("D11").Value = 'NA'
$ws.Range("E11").Value = 'NA'
$ws.Range("F11").Value = 'NA'
$ws.Range("G11").Value = 'F'

# Row 12
$ws.Range("C12").Value = 'BIO 114 or BOT 121.'
$ws.Range("D12").Value = 'NA'
$ws.Range("E12").Value = 'NA'
$ws.Range("F12").Value = 'NA'
$ws.Range("G12").Value = 'SP'

# Row 13
$ws.Range("C13").Value = 'LA 202, LA 243.'
$ws.Range("D13").Value = 'MATH 118.'
$ws.Range("E13").Value = 'LA 203.'
$ws.Range("F13").Value = 'NA'
$ws.Range("G13").Value = 'W  '

# Row 14
$ws.Range("C14").Value = 'LA 203, LA 241.'
$ws.Range("D14").Value = 'NA'
$ws.Range("E14").Value = 'LA 204.'
$ws.Range("F14").Value = 'NA'
$ws.Range("G14").Value = 'SP '

# Row 15
$ws.Range("D15").Value = 'LA 170.'
$ws.Range("E15").Value = 'LA 202.'
$ws.Range("F15").Value = 'NA'
$ws.Range("G15").Value = 'F '

# Row 16
$ws.Range("D16").Value = 'NA'
$ws.Range("E16").Value = 'NA'
$ws.Range("F16").Value = 'NA'
$ws.Range("G16").Value = 'TBD'

# Row 17
$ws.Range("D17").Value = 'NA'
$ws.Range("E17").Value = 'NA'
$ws.Range("F17").Value = 'NA'
$ws.Range("G17").Value = 'TBD'

# Row 18
$ws.Range("C18").Value = 'LA 211, LA 212, or consent of instructor.'
$ws.Range("D18").Value = 'NA'
$ws.Range("E18").Value = 'NA'
$ws.Range("F18").Value = 'NA'
$ws.Range("G18").Value = 'F'

# Row 19
$ws.Range("C19").Value = 'LA 211, LA 212.'
$ws.Range("D19").Value = 'NA'
$ws.Range("E19").Value = 'NA'
$ws.Range("F19").Value = 'NA'
$ws.Range("G19").Value = 'F'

# Row 20
$ws.Range("C20").Value = 'AEPS 233; AEPS 234; and AEPS 381 or LA 221.'
$ws.Range("D20").Value = 'NA'
$ws.Range("E20").Value = 'NA'
$ws.Range("F20").Value = 'NA'
$ws.Range("G20").Value = 'SP'

# Row 21
$ws.Range("C21").Value = 'LA 242.'
$ws.Range("D21").Value = 'NA'
$ws.Range("E21").Value = 'NA'
$ws.Range("F21").Value = 'NA'
$ws.Range("G21").Value = 'SP'

# Row 22
$ws.Range("C22").Value = 'LA 204.'
$ws.Range("D22").Value = 'NA'
$ws.Range("E22").Value = 'NA'
$ws.Range("F22").Value = 'NA'
$ws.Range("G22").Value = 'W'

# Row 23
$ws.Range("D23").Value = 'NA'
$ws.Range("E23").Value = 'NA'
$ws.Range("F23").Value = 'NA'
$ws.Range("G23").Value = 'F, W, SP'

# Row 24
$ws.Range("D24").Value = 'NA'
$ws.Range("E24").Value = 'NA'
$ws.Range("F24").Value = 'NA'
$ws.Range("G24").Value = 'F, W, SP'

# Row 25
$ws.Range("C25").Value = 'Completion of four design focus studios (16 units from LA 402 - LA 405).'
$ws.Range("D25").Value = 'NA'
$ws.Range("E25").Value = 'NA'
$ws.Range("F25").Value = 'NA'
$ws.Range("G25").Value = 'F, SP'

# Row 26
$ws.Range("C26").Value = 'LA 204, LA 211, LA 212 or consent of instructor.'
$ws.Range("D26").Value = 'NA'
$ws.Range("E26").Value = 'Integrated Learning Course (ILC) of student''s option.'
$ws.Range("F26").Value = 'NA'
$ws.Range("G26").Value = 'F, W, SU '

# Row 27
$ws.Range("C27").Value = 'LA 204, LA 211, LA 212, LA 402, and LA 405.'
$ws.Range("D27").Value = 'LA 220.'
$ws.Range("E27").Value = 'Integrated Learning Course (ILC) of student''s option.'
$ws.Range("F27").Value = 'NA'
$ws.Range("G27").Value = 'W, SP  '

# Row 28
$ws.Range("C28").Value = 'LA 204, LA 211, LA 212, LA 402, and LA 405.'
$ws.Range("D28").Value = 'NA'
$ws.Range("E28").Value = 'Integrated Learning Course (ILC) of student''s option.'
$ws.Range("F28").Value = 'NA'
$ws.Range("G28").Value = 'F, W, SP '

# Row 29
$ws.Range("C29").Value = 'LA 204, LA 243, LA 242, LA 241 or consent of instructor.'
$ws.Range("D29").Value = 'NA'
$ws.Range("E29").Value = 'Integrated Learning Course (ILC) of student''s option.'
$ws.Range("F29").Value = 'NA'
$ws.Range("G29").Value = 'F, W '

# Row 30
$ws.Range("C30").Value = 'LA 402, LA 405 or consent of instructor.'
$ws.Range("D30").Value = 'NA'
$ws.Range("E30").Value = 'Integrated Learning Course (ILC) of student''s option.'
$ws.Range("F30").Value = 'NA'
$ws.Range("G30").Value = 'W, SP '

# Row 31
$ws.Range("C31").Value = 'LA 170, LA 204 or consent of instructor.'
$ws.Range("D31").Value = 'NA'
$ws.Range("E31").Value = 'NA'
$ws.Range("F31").Value = 'NA'
$ws.Range("G31").Value = 'F, W, SP'

# Row 32
$ws.Range("D32").Value = 'LA 220.'
$ws.Range("E32").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F32").Value = 'NA'
$ws.Range("G32").Value = 'F, W, SP '

# Row 33
$ws.Range("C33").Value = 'LA 211.'
$ws.Range("D33").Value = 'NA'
$ws.Range("E33").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F33").Value = 'NA'
$ws.Range("G33").Value = 'F, W, SP '

# Row 34
$ws.Range("C34").Value = 'LA 241, LA 242, LA 243 or consent of instructor.'
$ws.Range("D34").Value = 'NA'
$ws.Range("E34").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F34").Value = 'NA'
$ws.Range("G34").Value = 'F, W, SP '

# Row 35
$ws.Range("D35").Value = 'NA'
$ws.Range("E35").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F35").Value = 'NA'
$ws.Range("G35").Value = 'F, W, SP'

# Row 36
$ws.Range("C36").Value = 'LA 170, LA 202, LA 203, LA 204 or consent of instructor.'
$ws.Range("D36").Value = 'NA'
$ws.Range("E36").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F36").Value = 'NA'
$ws.Range("G36").Value = 'F, W, SP '

# Row 37
$ws.Range("C37").Value = 'LA 170, LA 204 or consent of instructor.'
$ws.Range("D37").Value = 'NA'
$ws.Range("E37").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F37").Value = 'NA'
$ws.Range("G37").Value = 'F, W, SP '

# Row 38
$ws.Range("D38").Value = 'LA 220.'
$ws.Range("E38").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F38").Value = 'NA'
$ws.Range("G38").Value = 'F, W, SP '

# Row 39
$ws.Range("C39").Value = 'LA 221.'
$ws.Range("D39").Value = 'NA'
$ws.Range("E39").Value = 'Design Focus Studio of student''s option.'
$ws.Range("F39").Value = 'NA'
$ws.Range("G39").Value = 'F, W, SP '

# Row 40
$ws.Range("C40").Value = 'Completion of Design Focus Sequence (20 units from LA 402-LA 405).'
$ws.Range("D40").Value = 'NA'
$ws.Range("E40").Value = 'NA'
$ws.Range("F40").Value = 'NA'
$ws.Range("G40").Value = 'F, W, SP'

# Row 41
$ws.Range("D41").Value = 'NA'
$ws.Range("E41").Value = 'NA'
$ws.Range("F41").Value = 'NA'
$ws.Range("G41").Value = 'TBD'

# Row 42
$ws.Range("D42").Value = 'NA'
$ws.Range("E42").Value = 'NA'
$ws.Range("F42").Value = 'NA'
$ws.Range("G42").Value = 'TBD'

# Row 43
$ws.Range("D43").Value = 'NA'
$ws.Range("E43").Value = 'NA'
$ws.Range("F43").Value = 'NA'
$ws.Range("G43").Value = 'TBD'

# Row 44
$ws.Range("D44").Value = 'NA'
$ws.Range("E44").Value = 'NA'
$ws.Range("F44").Value = 'NA'
$ws.Range("G44").Value = 'TBD'

# Row 45
$ws.Range("D45").Value = 'NA'
$ws.Range("E45").Value = 'NA'
$ws.Range("F45").Value = 'NA'
$ws.Range("G45").Value = 'TBD'

# Row 46
$ws.Range("D46").Value = 'NA'
$ws.Range("E46").Value = 'NA'
$ws.Range("F46").Value = 'NA'
$ws.Range("G46").Value = 'TBD'

# Row 47
$ws.Range("D47").Value = 'NA'
$ws.Range("E47").Value = 'NA'
$ws.Range("F47").Value = 'NA'
$ws.Range("G47").Value = 'TBD'

# Row 48
$ws.Range("D48").Value = 'NA'
$ws.Range("E48").Value = 'NA'
$ws.Range("F48").Value = 'NA'
$ws.Range("G48").Value = 'TBD'
